# Edit script generated to reproduce the target OOXML diff for before.pptx
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Reposition existing shapes ---
$shp = $s.Shapes.Item(3)  # id=26
$shp.Top = 8.109055118110236
$shp = $s.Shapes.Item(4)  # id=15
$shp.Left = 57.06047439575195
$shp.Top = 144.66458129882812
$shp = $s.Shapes.Item(5)  # id=16
$shp.Top = 116.54701232910156
$shp = $s.Shapes.Item(6)  # id=47
$shp.Top = 280.0817565917969
$shp = $s.Shapes.Item(7)  # id=48
$shp.Top = 280.0817565917969
$shp = $s.Shapes.Item(8)  # id=49
$shp.Left = 644.6738891601562
$shp.Top = 280.0817565917969
$shp = $s.Shapes.Item(9)  # id=50
$shp.Top = 280.0817565917969
$shp = $s.Shapes.Item(11)  # id=53
$shp.Top = 37.53488540649414
$shp = $s.Shapes.Item(12)  # id=54
$shp.Top = 101.22582677165354
$shp = $s.Shapes.Item(13)  # id=57
$shp.Left = 17.218740157480315
$shp.Top = 258.98663330078125

# --- Add new shapes by copying the "fold" label style (id=24) and repositioning/retexting ---
$template = $s.Shapes.Item(1)
$newShapes = @()
$template.Copy()
$pasted = $s.Shapes.Paste()
$ns = $pasted.Item(1)
$ns.Left = 227.00103759765625
$ns.Top = 48.47803497314453
$ns.TextFrame.TextRange.Text = "hollow"
$newShapes += $ns
$template.Copy()
$pasted = $s.Shapes.Paste()
$ns = $pasted.Item(1)
$ns.Left = 437.42010498046875
$ns.Top = 48.47803497314453
$ns.TextFrame.TextRange.Text = "bridge"
$newShapes += $ns
$template.Copy()
$pasted = $s.Shapes.Paste()
$ns = $pasted.Item(1)
$ns.Left = 640.4968872070312
$ns.Top = 45.091181102362206
$ns.TextFrame.TextRange.Text = "atop"
$newShapes += $ns

# --- Add click-entrance animations ("Appear") for the new shapes, in order ---
$mainSeq = $s.TimeLine.MainSequence
foreach ($ns in $newShapes) {
    $mainSeq.AddEffect($ns, 1) | Out-Null
}

Write-Host "Edit complete"
